$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.080416947809462
$ws.Range("D2").Value = 1.082600923681428
$ws.Range("E2").Value = 1.083516436617712
$ws.Range("F2").Value = 1.093851892338386
$ws.Range("I2").Value = 1.06169137221202
$ws.Range("J2").Value = 1.085296588463696
$ws.Range("K2").Value = 1.085268716931678
$ws.Range("L2").Value = 1.086181849270659
$ws.Range("M2").Value = 1.096490727932794

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.081638820709448
$ws.Range("D3").Value = 1.083590410843402
$ws.Range("E3").Value = 1.084587261377797
$ws.Range("F3").Value = 1.094955990059693
$ws.Range("I3").Value = 1.062087799923796
$ws.Range("J3").Value = 1.086178520799056
$ws.Range("K3").Value = 1.086076453293011
$ws.Range("L3").Value = 1.087070899070885
$ws.Range("M3").Value = 1.097414893440484

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.08242949720248
$ws.Range("D4").Value = 1.084230661042608
$ws.Range("E4").Value = 1.085280408757778
$ws.Range("F4").Value = 1.095670700494448
$ws.Range("I4").Value = 1.062343137207977
$ws.Range("J4").Value = 1.086748650882023
$ws.Range("K4").Value = 1.086598477465925
$ws.Range("L4").Value = 1.087645820267207
$ws.Range("M4").Value = 1.098012574447813

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.082761909246758
$ws.Range("D5").Value = 1.084499819287961
$ws.Range("E5").Value = 1.08557186882606
$ws.Range("F5").Value = 1.095971233636998
$ws.Range("I5").Value = 1.062450199306736
$ws.Range("J5").Value = 1.086988205182483
$ws.Range("K5").Value = 1.086817784879367
$ws.Range("L5").Value = 1.087887433080355
$ws.Range("M5").Value = 1.098263764481661

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.082817723413897
$ws.Range("D6").Value = 1.084545011982561
$ws.Range("E6").Value = 1.085620809841619
$ws.Range("F6").Value = 1.096021698529139
$ws.Range("I6").Value = 1.062468158994482
$ws.Range("J6").Value = 1.087028419894739
$ws.Range("K6").Value = 1.086854598687778
$ws.Range("L6").Value = 1.087927996022097
$ws.Range("M6").Value = 1.098305936008283

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.082433938861021
$ws.Range("D7").Value = 1.084234257559231
$ws.Range("E7").Value = 1.085284303021456
$ws.Range("F7").Value = 1.095674715961671
$ws.Range("I7").Value = 1.06234456888346
$ws.Range("J7").Value = 1.086751852321927
$ws.Range("K7").Value = 1.086601408457126
$ws.Range("L7").Value = 1.08764904903902
$ws.Range("M7").Value = 1.098015931155575

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.08082987737841
$ws.Range("D8").Value = 1.082935328704265
$ws.Range("E8").Value = 1.083878274975155
$ws.Range("F8").Value = 1.094224968576909
$ws.Range("I8").Value = 1.061825590901966
$ws.Range("J8").Value = 1.085594753659825
$ws.Range("K8").Value = 1.085541827074231
$ws.Range("L8").Value = 1.086482381341756
$ws.Range("M8").Value = 1.096803119733042

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.078003585153308
$ws.Range("D9").Value = 1.080646321009909
$ws.Range("E9").Value = 1.081402580407304
$ws.Range("F9").Value = 1.091672488852449
$ws.Range("I9").Value = 1.060902049681847
$ws.Range("J9").Value = 1.083551638393818
$ws.Range("K9").Value = 1.083669819673396
$ws.Range("L9").Value = 1.084423831009779
$ws.Range("M9").Value = 1.09466354645615

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.076119484360374
$ws.Range("D10").Value = 1.079120197503143
$ws.Range("E10").Value = 1.079753361865137
$ws.Range("F10").Value = 1.089972250070596
$ws.Range("I10").Value = 1.060280259538646
$ws.Range("J10").Value = 1.082186719757437
$ws.Range("K10").Value = 1.082418491065267
$ws.Range("L10").Value = 1.083049588579778
$ws.Range("M10").Value = 1.093235482703579

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.075303650135021
$ws.Range("D11").Value = 1.07845933321326
$ws.Range("E11").Value = 1.079039518038701
$ws.Range("F11").Value = 1.089236353995435
$ws.Range("I11").Value = 1.060009566565047
$ws.Range("J11").Value = 1.081595010112787
$ws.Range("K11").Value = 1.081875854864655
$ws.Range("L11").Value = 1.082454072334561
$ws.Range("M11").Value = 1.092616707253953

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.075000610225723
$ws.Range("D12").Value = 1.07821385158801
$ws.Range("E12").Value = 1.078774405635406
$ws.Range("F12").Value = 1.088963056371649
$ws.Range("I12").Value = 1.059908800229038
$ws.Range("J12").Value = 1.081375117985726
$ws.Range("K12").Value = 1.081674174179077
$ws.Range("L12").Value = 1.082232801052846
$ws.Range("M12").Value = 1.092386803344554

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.075065613468904
$ws.Range("D13").Value = 1.078266508568451
$ws.Range("E13").Value = 1.078831271305858
$ws.Range("F13").Value = 1.089021677530015
$ws.Range("I13").Value = 1.059930424864547
$ws.Range("J13").Value = 1.081422290353154
$ws.Range("K13").Value = 1.081717440879809
$ws.Range("L13").Value = 1.082280267677366
$ws.Range("M13").Value = 1.092436121385875

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.075278600828319
$ws.Range("D14").Value = 1.078439041777001
$ws.Range("E14").Value = 1.079017602933968
$ws.Range("F14").Value = 1.089213762176226
$ws.Range("I14").Value = 1.060001241656954
$ws.Range("J14").Value = 1.081576835905168
$ws.Range("K14").Value = 1.08185918635852
$ws.Range("L14").Value = 1.082435783415199
$ws.Range("M14").Value = 1.092597704618627

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.075409828950438
$ws.Range("D15").Value = 1.07854534420324
$ws.Range("E15").Value = 1.079132413403036
$ws.Range("F15").Value = 1.089332118059628
$ws.Range("I15").Value = 1.06004484519682
$ws.Range("J15").Value = 1.081672042600232
$ws.Range("K15").Value = 1.081946504308779
$ws.Range("L15").Value = 1.082531592487219
$ws.Range("M15").Value = 1.092697252985656

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.076173628305318
$ws.Range("D16").Value = 1.079164055922363
$ws.Range("E16").Value = 1.079800743124408
$ws.Range("F16").Value = 1.090021095702848
$ws.Range("I16").Value = 1.060298193868043
$ws.Range("J16").Value = 1.082225974934593
$ws.Range("K16").Value = 1.082454487059195
$ws.Range("L16").Value = 1.083089101259543
$ws.Range("M16").Value = 1.093276539984573

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.076652736523996
$ws.Range("D17").Value = 1.079552145443571
$ws.Range("E17").Value = 1.080220042684818
$ws.Range("F17").Value = 1.090453357813935
$ws.Range("I17").Value = 1.060456723232164
$ws.Range("J17").Value = 1.082573256365616
$ws.Range("K17").Value = 1.08277291572442
$ws.Range("L17").Value = 1.083438687943874
$ws.Range("M17").Value = 1.093639799831212

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.076932191771472
$ws.Range("D18").Value = 1.07977850742401
$ws.Range("E18").Value = 1.080464639973408
$ws.Range("F18").Value = 1.09070551978553
$ws.Range("I18").Value = 1.060549050497297
$ws.Range("J18").Value = 1.082775753158741
$ws.Range("K18").Value = 1.082958572389828
$ws.Range("L18").Value = 1.083642551599134
$ws.Range("M18").Value = 1.09385164322829

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.077027478764915
$ws.Range("D19").Value = 1.079855690342809
$ws.Range("E19").Value = 1.080548045925973
$ws.Range("F19").Value = 1.090791505752575
$ws.Range("I19").Value = 1.060580507938558
$ws.Range("J19").Value = 1.082844788068533
$ws.Range("K19").Value = 1.083021863414497
$ws.Range("L19").Value = 1.083712056343035
$ws.Range("M19").Value = 1.093923869607293

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.076601332808011
$ws.Range("D20").Value = 1.079510507537382
$ws.Range("E20").Value = 1.080175053028997
$ws.Range("F20").Value = 1.090406977001387
$ws.Range("I20").Value = 1.060439729039994
$ws.Range("J20").Value = 1.08253600324965
$ws.Range("K20").Value = 1.082738759352766
$ws.Range("L20").Value = 1.083401185194002
$ws.Range("M20").Value = 1.093600829623505

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.075215881450223
$ws.Range("D21").Value = 1.078388235252481
$ws.Range("E21").Value = 1.078962731789661
$ws.Range("F21").Value = 1.089157196746938
$ws.Range("I21").Value = 1.059980393923774
$ws.Range("J21").Value = 1.081531328993791
$ws.Range("K21").Value = 1.08181744920711
$ws.Range("L21").Value = 1.082389989846952
$ws.Range("M21").Value = 1.092550124134547

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.074344776360514
$ws.Range("D22").Value = 1.07768257606207
$ws.Range("E22").Value = 1.07820073395221
$ws.Range("F22").Value = 1.088371681450418
$ws.Range("I22").Value = 1.059690324892975
$ws.Range("J22").Value = 1.080899043455633
$ws.Range("K22").Value = 1.081237481918057
$ws.Range("L22").Value = 1.081753805567564
$ws.Range("M22").Value = 1.091889138044788

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.074806567856661
$ws.Range("D23").Value = 1.078056663598815
$ws.Range("E23").Value = 1.078604661290391
$ws.Range("F23").Value = 1.088788072574326
$ws.Range("I23").Value = 1.059844216222074
$ws.Range("J23").Value = 1.081234287857264
$ws.Range("K23").Value = 1.081545000474659
$ws.Range("L23").Value = 1.082091097638278
$ws.Range("M23").Value = 1.092239574254061

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.076624559935571
$ws.Range("D24").Value = 1.079529321927202
$ws.Range("E24").Value = 1.080195381834328
$ws.Range("F24").Value = 1.09042793439867
$ws.Range("I24").Value = 1.060447408417046
$ws.Range("J24").Value = 1.082552836534639
$ws.Range("K24").Value = 1.082754193385702
$ws.Range("L24").Value = 1.083418131209095
$ws.Range("M24").Value = 1.093618438705975

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.078734226506662
$ws.Range("D25").Value = 1.081238102912341
$ws.Range("E25").Value = 1.082042384705401
$ws.Range("F25").Value = 1.092332113952471
$ws.Range("I25").Value = 1.061141880158449
$ws.Range("J25").Value = 1.084080330242828
$ws.Range("K25").Value = 1.08415436162821
$ws.Range("L25").Value = 1.084956343393022
$ws.Range("M25").Value = 1.095216971038559
